$wb = $excel.ActiveWorkbook

# --- Fase de Grupos: enter the results for the last two Group B matches ---
$wsGrupos = $wb.Worksheets.Item("Fase de Grupos")
$wsGrupos.Unprotect("CC01")

$wsGrupos.Range("F13").Value = 0
$wsGrupos.Range("H13").Value = 2
$wsGrupos.Range("F14").Value = 1
$wsGrupos.Range("H14").Value = 0

$wsGrupos.Protect("CC01")

# --- Finais: fix the broken Round-of-16 bracket reference for 2nd place Group B ---
$wsFinais = $wb.Worksheets.Item("Finais")
$wsFinais.Unprotect("CC01")

$wsFinais.Range("K6").Formula = "=IF('Fase de Grupos'!AH12=3,'Fase de Grupos'!AF12,""2º do Grupo B"")"

$wsFinais.Protect("CC01")

$excel.Calculate()
